$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where course_duration (C) becomes "course duration not found"
$durationNotFoundRows = 2,3,4,5,6,8,10,11,12,13
foreach ($r in $durationNotFoundRows) {
    $ws.Cells.Item($r, 3).Value = "course duration not found"
}

# Row 7: course_duration becomes "Course Objective"
$ws.Cells.Item(7, 3).Value = "Course Objective"

# course_summary (D) updates -- force text storage for numeric-looking values
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "1"

$ws.Cells.Item(3, 4).Value = "course summary not found"
$ws.Cells.Item(4, 4).Value = "course summary not found"
$ws.Cells.Item(5, 4).Value = "course summary not found"
$ws.Cells.Item(6, 4).Value = "course summary not found"

$ws.Cells.Item(7, 4).Value = "Lead Developer Participants should gain competancy in using core techniques to handle natural language content to undertake analysis to detect patterns and derive insights for development of applications like mentioned in course summary Course Objective Natural Language Processing is a sub-field of Artificial Intelligence. It is used for processing and analysing large amounts of natural language. Some applications include search engines (Google), text classification (spam filters), identifying sentiments for a product (sentiment analysis), methods for discovering abstract topics in a collection of documents (topic modelling) and machine translation technologies. Concepts covered include cleaning, exploring datasets through methods rooted in Corpus Linguistics, and application of feature engineering techniques to transform textual data into a numerical representation. Key techniques such as word embeddings and language modelling are also introduced as well as illustrations as to how they can be performed over a dataset."

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "1.1"

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "1"

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "1"

$ws.Cells.Item(12, 4).Value = "course summary not found"
$ws.Cells.Item(13, 4).Value = "course summary not found"
